$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: participant renamed from Hebrew "שחקן ג" to "Sahkan" (Latin text loses the right-alignment style)
$ws.Range("C4").Value = "Sahkan"
$ws.Range("C4").HorizontalAlignment = 1   # xlHAlignGeneral -> style s=1

# Rows 5-19: downstream participant labels renumber after the C/D merge + new "Player Y" (Row10) insertion
$ws.Range("C5").Value = "שחקן ה"
$ws.Range("C6").Value = "שחקן ו"
$ws.Range("C7").Value = "שחקן ז"
$ws.Range("C8").Value = "שחקן ח"
$ws.Range("C9").Value = "שחקן ט"
$ws.Range("C10").Value = "שחקן י"
$ws.Range("C11").Value = "שחקן יא"
$ws.Range("C12").Value = "שחקן יב"
$ws.Range("C13").Value = "שחקן יג"
$ws.Range("C14").Value = "שחקן יד"
$ws.Range("C15").Value = "שחקן טו"
$ws.Range("C16").Value = "שחקן טז"
$ws.Range("C17").Value = "שחקן יז"
$ws.Range("C18").Value = "שחקן יח"
$ws.Range("C19").Value = "שחקן יט"

# Rows 20-31: newly-entered participant names (already right-aligned style s=3)
$ws.Range("C20").Value = "שחקן כ"
$ws.Range("C21").Value = "שחקן כא"
$ws.Range("C22").Value = "שחקן כב"
$ws.Range("C23").Value = "שחקן כג"
$ws.Range("C24").Value = "שחקן כד"
$ws.Range("C25").Value = "שחקן כה"
$ws.Range("C26").Value = "שחקן כו"
$ws.Range("C27").Value = "שחקן כז"
$ws.Range("C28").Value = "שחקן כח"
$ws.Range("C29").Value = "שחקן כט"
$ws.Range("C30").Value = "שחקן ל"
$ws.Range("C31").Value = "שחקן לא"

# Rows 32-50: newly-entered participant names; alignment switches to right-aligned (style s=3)
$ws.Range("C32").Value = "שחקן לב"
$ws.Range("C32").HorizontalAlignment = -4152   # xlHAlignRight -> style s=3
$ws.Range("C33").Value = "שחקן לג"
$ws.Range("C33").HorizontalAlignment = -4152   # xlHAlignRight -> style s=3
$ws.Range("C34").Value = "שחקן לד"
$ws.Range("C34").HorizontalAlignment = -4152   # xlHAlignRight -> style s=3
$ws.Range("C35").Value = "שחקן לה"
$ws.Range("C35").HorizontalAlignment = -4152   # xlHAlignRight -> style s=3
$ws.Range("C36").Value = "שחקן לו"
$ws.Range("C36").HorizontalAlignment = -4152   # xlHAlignRight -> style s=3
$ws.Range("C37").Value = "שחקן לז"
$ws.Range("C37").HorizontalAlignment = -4152   # xlHAlignRight -> style s=3
$ws.Range("C38").Value = "שחקן לח"
$ws.Range("C38").HorizontalAlignment = -4152   # xlHAlignRight -> style s=3
$ws.Range("C39").Value = "שחקן לט"
$ws.Range("C39").HorizontalAlignment = -4152   # xlHAlignRight -> style s=3
$ws.Range("C40").Value = "שחקן מ"
$ws.Range("C40").HorizontalAlignment = -4152   # xlHAlignRight -> style s=3
$ws.Range("C41").Value = "שחקן מא"
$ws.Range("C41").HorizontalAlignment = -4152   # xlHAlignRight -> style s=3
$ws.Range("C42").Value = "שחקן מב"
$ws.Range("C42").HorizontalAlignment = -4152   # xlHAlignRight -> style s=3
$ws.Range("C43").Value = "שחקן מג"
$ws.Range("C43").HorizontalAlignment = -4152   # xlHAlignRight -> style s=3
$ws.Range("C44").Value = "שחקן מד"
$ws.Range("C44").HorizontalAlignment = -4152   # xlHAlignRight -> style s=3
$ws.Range("C45").Value = "שחקן מה"
$ws.Range("C45").HorizontalAlignment = -4152   # xlHAlignRight -> style s=3
$ws.Range("C46").Value = "שחקן מו"
$ws.Range("C46").HorizontalAlignment = -4152   # xlHAlignRight -> style s=3
$ws.Range("C47").Value = "שחקן מז"
$ws.Range("C47").HorizontalAlignment = -4152   # xlHAlignRight -> style s=3
$ws.Range("C48").Value = "שחקן מח"
$ws.Range("C48").HorizontalAlignment = -4152   # xlHAlignRight -> style s=3
$ws.Range("C49").Value = "שחקן מט"
$ws.Range("C49").HorizontalAlignment = -4152   # xlHAlignRight -> style s=3
$ws.Range("C50").Value = "שחקן נ"
$ws.Range("C50").HorizontalAlignment = -4152   # xlHAlignRight -> style s=3

# Restore portrait page orientation (adds <pageSetup orientation="portrait".../>)
$ws.PageSetup.Orientation = 1   # xlPortrait

# Move the active selection/view: frozen pane resets to A2, selection to D47
$ws.Range("D47").Select() | Out-Null
